$wb = $excel.ActiveWorkbook

# --- "Weekly Sales" sheet: append a new weekly row (row 20) ---
$ws1 = $wb.Worksheets.Item("Weekly Sales")
$ws1.Range("A20").Value = 45662.99999999999
$ws1.Range("A20").NumberFormat = $ws1.Range("A19").NumberFormat
$ws1.Range("B20").Value = 3

# --- "Merged (Optional)" sheet: append a new matching row (row 28) ---
$ws3 = $wb.Worksheets.Item("Merged (Optional)")
$ws3.Range("A28").Value = 45662.99999999999
$ws3.Range("A28").NumberFormat = $ws3.Range("A27").NumberFormat
$ws3.Range("B28").Value = 3
$ws3.Range("C28").Value = 0
